# "Error Calculations and Plots"
#
# The source data table (ID, A, B, C, D, F) contained two rows that were
# removed from the cleaned/missing-data export:
#   - row 26 "RM 232"
#   - row 28 "SC 92"
# Removing them shifts every following row up (first by one, then by one
# more), and two previously-missing values get filled in as part of the
# same edit:
#   - column B for "SC 5"   (ends up on row 26) becomes -20.2
#   - column F for "SC 232" (ends up on row 33) becomes 17.53
# "SC 101" (ends up on row 27) loses its column-B value, becoming missing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete "RM 232" (row 26). Everything below shifts up by one row, so the
# row that used to be 28 ("SC 92") is now row 27.
$ws.Rows(26).Delete()

# Delete "SC 92", now sitting on row 27. Everything below shifts up again.
$ws.Rows(27).Delete()

# Fill in the previously-missing value for "SC 5" (now row 26), column B.
$ws.Range("B26").Value = -20.2

# "SC 101" (now row 27) loses its column-B value -- clear it back to missing.
$ws.Range("B27").ClearContents()

# Fill in the previously-missing value for "SC 232" (now row 33), column F.
$ws.Range("F33").Value = 17.53
